# Correção e atualização 13/11
# Insert 8 new rows (L2L3 / MUX 25 / channels 1-8) above the existing "L4" row
# (previously row 5, now shifted down to row 13), and update the active cell
# selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 8 new blank rows at rows 5:12 (pushes the old rows 5-29 down to 13-37).
# The default Insert() behaviour copies per-column formatting from the row
# directly above (row 4), which already matches the styles used by the rest
# of the data rows (center/middle for columns A, B, D and center-only for C).
$ws.Rows("5:12").Insert()

# Fill in the new rows with the "L2L3" location data (no "Nivel" / column B value).
for ($i = 0; $i -lt 8; $i++) {
    $row = 5 + $i
    $ws.Range("A$row").Value = "L2L3"
    $ws.Range("B$row").Clear()
    $ws.Range("C$row").Value = 25
    $ws.Range("D$row").Value = $i + 1
}

# Update the selected cell, matching the saved selection in the workbook.
$null = $ws.Range("G14").Select()
